$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as literal text instead of
# auto-converting numeric-looking strings (dropping trailing zeros, mangling
# strings that contain subscript digits, etc.)

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'64.589.72"
$ws.Range("E2").Value = "  -5.43%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.355.33"
$ws.Range("E3").Value = "  -7.59%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'185.75"
$ws.Range("E5").Value = "  -8.55%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'527.44"
$ws.Range("E6").Value = "  -8.64%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.596"
$ws.Range("E7").Value = "  -3.92%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "'3.346.77"
$ws.Range("E8").Value = "  -7.69%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "'0.621"
$ws.Range("E10").Value = "  -9.69%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'58.37"
$ws.Range("E11").Value = "  -4.93%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.132"
$ws.Range("E12").Value = "  -11.61%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -11.17%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'9.17"
$ws.Range("E14").Value = "  -11.57%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.878.28"
$ws.Range("E15").Value = "  -7.71%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -4.53%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'3.341.71"
$ws.Range("E17").Value = "  -7.89%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'64.328.56"
$ws.Range("E18").Value = "  -5.55%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'17.26"
$ws.Range("E19").Value = "  -11.23%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'10.98"
$ws.Range("E20").Value = "  -11.64%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "'0.957"
$ws.Range("E21").Value = "  -11.17%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'373.96"
$ws.Range("E22").Value = "  -8.90%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "'3.72"
$ws.Range("E23").Value = "  -12.11%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'80.41"
$ws.Range("E24").Value = "  -6.21%  "

# Row 25 - RenderToken
$ws.Range("D25").Value = "'10.85"
$ws.Range("E25").Value = "  -17.81%  "

# Row 26 - was LEO, now Toncoin
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'3.80"
$ws.Range("E26").Value = "  -4.86%  "

# Row 27 - was Toncoin, now LEO
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'6.02"
$ws.Range("E27").Value = "  -1.87%  "

# Row 28 - ImmutableX
$ws.Range("D28").Value = "'2.64"
$ws.Range("E28").Value = "  -10.09%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'11.27"
$ws.Range("E29").Value = "  -10.89%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'8.37"
$ws.Range("E30").Value = "  -10.76%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'28.83"
$ws.Range("E31").Value = "  -9.45%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'652.41"
$ws.Range("E32").Value = "  -4.40%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "'6.75"
$ws.Range("E33").Value = "  -13.47%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "'11.10"
$ws.Range("E34").Value = "  -9.45%  "

# Row 35 - OKB
$ws.Range("D35").Value = "'59.64"
$ws.Range("E35").Value = "  -6.66%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "'0.104"
$ws.Range("E36").Value = "  -9.61%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.01%  "

# Row 38 - InjectiveProtocol
$ws.Range("D38").Value = "'36.24"
$ws.Range("E38").Value = "  -13.74%  "

# Row 39 - TheGraph
$ws.Range("D39").Value = "'0.377"
$ws.Range("E39").Value = "  -9.49%  "

# Row 40 - FirstDigitalUSD
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41 - Kaspa
$ws.Range("D41").Value = "'0.124"
$ws.Range("E41").Value = "  -8.60%  "

# Row 42 - ThetaToken
$ws.Range("D42").Value = "'2.74"
$ws.Range("E42").Value = "  -14.64%  "

# Row 43 - Maker
$ws.Range("D43").Value = "'2.791.98"
$ws.Range("E43").Value = "  -12.49%  "

# Row 44 - PEPE
$ws.Range("D44").Value = "'0.0" + [char]0x2083 + "0623"
$ws.Range("E44").Value = "  -19.21%  "

# Row 45 - was WEMIXToken, now VeChain
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0389"
$ws.Range("E45").Value = "  -6.79%  "

# Row 46 - was VeChain, now WEMIXToken
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.59"
$ws.Range("E46").Value = "  -9.80%  "

# Row 47 - Fetch.AI
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = "  -14.67%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -6.27%  "

# Row 49 - Monero
$ws.Range("D49").Value = "'134.85"
$ws.Range("E49").Value = "  -3.34%  "

# Row 50 - Stacks
$ws.Range("D50").Value = "'2.64"
$ws.Range("E50").Value = "  -3.50%  "

# Row 51 - was ApeXProtocol, now dogwifhat
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.36"
$ws.Range("E51").Value = "  -17.95%  "
